# Fruta / hortaliza, semanal
# Insert this week's two new price rows (Española + Madrigal) at the top of
# the data block (row 10), pushing the historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 10 (old row 10 -> row 12, etc.)
$ws.Range("A10:A11").EntireRow.Insert()

# New row 10: Alcachofa, Española
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 45114
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100112013
$ws.Cells.Item(10, 7).Value = "Alcachofa"
$ws.Cells.Item(10, 8).Value = "Española"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 14000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 14500
$ws.Cells.Item(10, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 483
$ws.Cells.Item(10, 17).Value = 30
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# New row 11: Alcachofa, Madrigal
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 45114
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112013
$ws.Cells.Item(11, 7).Value = "Alcachofa"
$ws.Cells.Item(11, 8).Value = "Madrigal"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 400
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 16000
$ws.Cells.Item(11, 13).Value = 15500
$ws.Cells.Item(11, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 388
$ws.Cells.Item(11, 17).Value = 40
$ws.Cells.Item(11, 18).Value = "Hortaliza"
